# "Generate Report for Handoff" — localization-status report refresh.
#
# For the four "Ready for handoff" rows (46ed3f22*, a1bd38a2*, b3a1ff7c*,
# d6ca633e*) in both the zh-cn and de-de sheets:
#   - Priority (col E) flips from "low" to "ht"
#   - Latest Handoff Datetime (col H) is bumped to the new handoff timestamp
#     (the Overview sheet's "Latest HO Xliff Generate Date" column shares the
#     same string as de-de's handoff datetime, so it updates automatically).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = @(4, 5, 6, 7)

foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-24 22:31:56"

    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-24 22:32:02"

    # Overview's "Latest HO Xliff Generate Date" tracks the de-de handoff
    # datetime (they share the same underlying text), so refresh it too.
    $wsOverview.Range("G$r").Value = "2016-08-24 22:32:02"
}
